{"js": "// The diary document has four paragraphs:\n//   1. \"2022-6-24 \u661f\u671f\u4e94\"\n//   2. \"\u6674\"\n//   3. \"2022-6-25 \u661f\u671f\u516d\"\n//   4. \"\u96e8\"\n// The edit inserts two new paragraphs (\"\u96e8\" then \"\u661f\u671f\u5929\") right after the\n// \"2022-6-25 \u661f\u671f\u516d\" paragraph (i.e. before the old last paragraph), and\n// changes the text of the old last paragraph from \"\u96e8\" to \"\u661f\u671f\u4e00\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Find the paragraph that currently reads \"\u661f\u671f\u516d\" (the weekday line for\n// 2022-6-25) \u2014 the two new paragraphs get inserted right after it \u2014 and the\n// paragraph that currently reads \"\u96e8\" (the last paragraph), whose text gets\n// replaced with \"\u661f\u671f\u4e00\".\nlet weekdayParagraph = null;\nlet lastRainParagraph = null;\nfor (const p of paragraphs.items) {\n  const t = p.text.trim();\n  if (t.endsWith(\"\u661f\u671f\u516d\")) {\n    weekdayParagraph = p;\n  }\n  if (t === \"\u96e8\") {\n    lastRainParagraph = p;\n  }\n}\n\nif (!weekdayParagraph) {\n  throw new Error(\"Could not find the '\u661f\u671f\u516d' paragraph to insert after.\");\n}\nif (!lastRainParagraph) {\n  throw new Error(\"Could not find the '\u96e8' paragraph to edit.\");\n}\n\n// Insert the two new paragraphs after the \"\u661f\u671f\u516d\" paragraph, in order, so\n// the new paragraph marks do not inherit the \"\u96e8\" paragraph's formatting\n// (mirrors appending after the preceding paragraph rather than splitting\n// the old last paragraph).\nconst rainInserted = weekdayParagraph.insertParagraph(\"\u96e8\", Word.InsertLocation.after);\nrainInserted.insertParagraph(\"\u661f\u671f\u5929\", Word.InsertLocation.after);\n\n// Change the original paragraph's text to \"\u661f\u671f\u4e00\".\nlastRainParagraph.insertText(\"\u661f\u671f\u4e00\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The diary document has four paragraphs:\n#   1. \"2022-6-24 \u661f\u671f\u4e94\"\n#   2. \"\u6674\"\n#   3. \"2022-6-25 \u661f\u671f\u516d\"\n#   4. \"\u96e8\"\n# The edit inserts two new paragraphs (\"\u96e8\" then \"\u661f\u671f\u5929\") right after the\n# \"2022-6-25 \u661f\u671f\u516d\" paragraph, and changes the text of the old last\n# paragraph from \"\u96e8\" to \"\u661f\u671f\u4e00\".\n\n$d = $word.ActiveDocument\n\n# Locate the \"\u661f\u671f\u516d\" paragraph (insert point) and the \"\u96e8\" paragraph\n# (the one whose text becomes \"\u661f\u671f\u4e00\") by scanning all paragraphs, so the\n# script does not depend on a hard-coded paragraph index.\n$weekdayIndex = $null\n$rainIndex = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t.EndsWith(\"\u661f\u671f\u516d\")) {\n        $weekdayIndex = $i\n    }\n    if ($t -eq \"\u96e8\") {\n        $rainIndex = $i\n    }\n}\n\nif (-not $weekdayIndex) {\n    throw \"Could not find the '\u661f\u671f\u516d' paragraph to insert after.\"\n}\nif (-not $rainIndex) {\n    throw \"Could not find the '\u96e8' paragraph to edit.\"\n}\n\n# Insert the two new paragraphs right after the \"\u661f\u671f\u516d\" paragraph, in\n# order. Appending after the preceding paragraph (rather than splitting the\n# old last paragraph) keeps the new paragraph marks free of the old last\n# paragraph's formatting.\n$weekdayRange = $d.Paragraphs.Item($weekdayIndex).Range\n$weekdayRange.InsertParagraphAfter()\n$d.Paragraphs.Item($weekdayIndex + 1).Range.Text = \"\u96e8\"\n\n$rainRange = $d.Paragraphs.Item($weekdayIndex + 1).Range\n$rainRange.InsertParagraphAfter()\n$d.Paragraphs.Item($weekdayIndex + 2).Range.Text = \"\u661f\u671f\u5929\"\n\n# The paragraph that used to hold \"\u96e8\" has shifted down by two positions;\n# update its text to \"\u661f\u671f\u4e00\".\n$lastPara = $d.Paragraphs.Item($rainIndex + 2)\n$lastPara.Range.Text = \"\u661f\u671f\u4e00\"\n"}
